$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: N4:Q4 "ESFUERZO IDEAL" (Sprint 2) updated ---
$ws.Range("N4").Value = 240
$ws.Range("O4").Value = 180
$ws.Range("P4").Value = 120
$ws.Range("Q4").Value = 60

# --- Row 6: "Nº TAREAS" counts filled in for Sprint 1 (F6:L6) and Sprint 2 (N6:U6) ---
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 1

$ws.Range("N6").Value = 3
$ws.Range("O6").Value = 3
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 2

# --- Tasks rows 8-14: Estado changes from "Pendiente" to "Terminado" ---
$ws.Range("C8").Value = "Terminado"
$ws.Range("C9").Value = "Terminado"
$ws.Range("C10").Value = "Terminado"
$ws.Range("C11").Value = "Terminado"
$ws.Range("C12").Value = "Terminado"
$ws.Range("C13").Value = "Terminado"
$ws.Range("C14").Value = "Terminado"

# --- Row 8 (Tarea 1): L8 now tracked with 0, Sprint2 daily effort (P8:S8) cleared ---
$ws.Range("L8").Value = 0
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()

# --- Rows 10-13 (Tareas 2.1, 2.2, 2.3, 3): L column now tracked with 0 ---
$ws.Range("L10").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("L13").Value = 0

# --- Row 14 (Tarea 6): Responsable changed, Sprint2 daily effort reduced ---
$ws.Range("M14").Value = "Carlos, Alex, Oscar"
$ws.Range("N14").Value = 60
$ws.Range("O14").Value = 60
$ws.Range("P14").Value = 60
$ws.Range("Q14").Value = 60
$ws.Range("R14").Value = 30
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0

# --- Row 15 (Tarea 5): Responsable changed, Sprint2 daily effort increased (task not finished) ---
$ws.Range("M15").Value = "Carlos, Alex, Oscar"
$ws.Range("R15").Value = 120
$ws.Range("S15").Value = 120
$ws.Range("T15").Value = 120
$ws.Range("U15").Value = 120

# --- Row 16 (Tarea 4): new activity data added across K-U ---
$ws.Range("K16").Value = 60
$ws.Range("L16").Value = 60
$ws.Range("M16").Value = "Carlos"
$ws.Range("N16").Value = 60
$ws.Range("O16").Value = 60
$ws.Range("P16").Value = 60
$ws.Range("Q16").Value = 60
$ws.Range("R16").Value = 60
$ws.Range("S16").Value = 60
$ws.Range("T16").Value = 60
$ws.Range("U16").Value = 60

# --- Selection moved ---
$ws.Range("V19").Select()
